$d = $word.ActiveDocument

# Word's company name for "YASSAKA" is being spelled out in full, and the
# job title now starts on its own line after the company name.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$result = $find.Execute("YASSAKA – Analista", $true, $false, $false, $false, $false, $false, 1, $false, "Yassaka Consultoria e Treinamento em Gestão Empresarial  – ^lAnalista", 2)
Write-Host "Replace result: $result"

# Word marks the last edited location with a hidden "_GoBack" bookmark.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("Analista Comercial | Maio/2025", $true, $false, $false, $false, $false, $false, 1, $false, "", 0) | Out-Null
$editRange = $find2.Parent
$editRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $editRange) | Out-Null
Write-Host "Bookmark added"
